# Reorder the player rows on the active sheet.
# The underlying data set (player, position, team) is unchanged; only the
# row order for rows 2-15 changes (rows 16-19 already stay in place).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) rows 2-15 as [Player, Position, Team] triples.
$data = @{}
for ($r = 2; $r -le 15; $r++) {
    $data[$r] = @(
        $ws.Cells.Item($r, 1).Value2,
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2
    )
}

# New order for rows 2-15, expressed as the source row (from the table above)
# that should land in each destination row.
$order = @{
    2  = 12
    3  = 14
    4  = 15
    5  = 8
    6  = 3
    7  = 4
    8  = 5
    9  = 2
    10 = 6
    11 = 7
    12 = 9
    13 = 10
    14 = 11
    15 = 13
}

foreach ($destRow in 2..15) {
    $srcRow = $order[$destRow]
    $vals = $data[$srcRow]
    $ws.Cells.Item($destRow, 1).Value = $vals[0]
    $ws.Cells.Item($destRow, 2).Value = $vals[1]
    $ws.Cells.Item($destRow, 3).Value = $vals[2]
}
